$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 14705.8125
$ws.Range("I80").Value = 13097.5
$ws.Range("K80").Value = 39292.5
$ws.Range("M80").Value = -38294.5
$ws.Range("H83").Value = 14705.8125
$ws.Range("I83").Value = 13097.5
$ws.Range("K83").Value = 117877.5
$ws.Range("M83").Value = -112885.5
$ws.Range("H127").Value = 982.4375
$ws.Range("I127").Value = 711.2308
$ws.Range("J127").Value = 2157.6667
$ws.Range("K127").Value = 2133.6924
$ws.Range("L127").Value = 6473.000100000001
$ws.Range("M127").Value = 2826.3076
$ws.Range("N127").Value = -16393.0001
$ws.Range("H132").Value = 2992403
$ws.Range("I132").Value = 3264135.2
$ws.Range("K132").Value = 9792405.600000001
$ws.Range("M132").Value = -9789875.600000001
$ws.Range("H137").Value = 15600.131
$ws.Range("I137").Value = 31894.5
$ws.Range("J137").Value = 3066
$ws.Range("K137").Value = 95683.5
$ws.Range("L137").Value = 9198
$ws.Range("M137").Value = -93133.5
$ws.Range("N137").Value = -14298
$ws.Range("H141").Value = 1150
$ws.Range("I141").Value = 1150
$ws.Range("K141").Value = 3450
$ws.Range("M141").Value = 1730

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16952.85
$ws.Range("I32").Value = 16771.639
$ws.Range("K32").Value = 16771.639
$ws.Range("M32").Value = -16484.639
$ws.Range("H45").Value = 2812.175
$ws.Range("I45").Value = 1878.5
$ws.Range("K45").Value = 1878.5
$ws.Range("M45").Value = -1501.5
$ws.Range("H74").Value = 6000012
$ws.Range("I74").Value = 6000012
$ws.Range("K74").Value = 6000012
$ws.Range("M74").Value = -5999138
$ws.Range("H77").Value = 6000012
$ws.Range("I77").Value = 6000012
$ws.Range("K77").Value = 30000060
$ws.Range("M77").Value = -29995692
$ws.Range("H102").Value = 2319
$ws.Range("I102").Value = 1872.4
$ws.Range("K102").Value = 1872.4
$ws.Range("M102").Value = -250.4000000000001
$ws.Range("H106").Value = 44999.5
$ws.Range("J106").Value = 44999.5
$ws.Range("L106").Value = 44999.5
$ws.Range("N106").Value = -47523.5
$ws.Range("H122").Value = 2060.9795
$ws.Range("I122").Value = 2014.9783
$ws.Range("K122").Value = 6044.9349
$ws.Range("M122").Value = -3594.9349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1413.1428
$ws.Range("I5").Value = 378.4
$ws.Range("J5").Value = 4000
$ws.Range("K5").Value = 378.4
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = -265.4
$ws.Range("N5").Value = -4226
$ws.Range("H99").Value = 1337.7273
$ws.Range("I99").Value = 1238.8
$ws.Range("J99").Value = 1420.1666
$ws.Range("K99").Value = 1238.8
$ws.Range("L99").Value = 1420.1666
$ws.Range("M99").Value = 259.2
$ws.Range("N99").Value = -4416.1666
$ws.Range("H105").Value = 4699
$ws.Range("I105").Value = 3298.5
$ws.Range("K105").Value = 3298.5
$ws.Range("M105").Value = -1551.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 1465.3182
$ws.Range("J58").Value = 2811.4
$ws.Range("K58").Value = 1465.3182
$ws.Range("L58").Value = 2811.4
$ws.Range("M58").Value = -1262.3182
$ws.Range("N58").Value = -3217.4
$ws.Range("I136").Value = 1465.3182
$ws.Range("J136").Value = 2811.4
$ws.Range("K136").Value = 4395.9546
$ws.Range("L136").Value = 8434.200000000001
$ws.Range("M136").Value = -1845.9546
$ws.Range("N136").Value = -13534.2
$ws.Range("H138").Value = 85283.19
$ws.Range("J138").Value = 85283.19
$ws.Range("L138").Value = 85283.19
$ws.Range("N138").Value = -95563.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 32855962
$ws.Range("I4").Value = 42214160
$ws.Range("J4").Value = 11260126
$ws.Range("K4").Value = 126642480
$ws.Range("L4").Value = 33780378
$ws.Range("M4").Value = -126642368
$ws.Range("N4").Value = -33780602
$ws.Range("H137").Value = 2715.35
$ws.Range("J137").Value = 2118.3
$ws.Range("L137").Value = 6354.900000000001
$ws.Range("N137").Value = -16554.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 533.5
$ws.Range("I2").Value = 290.8
$ws.Range("J2").Value = 938
$ws.Range("K2").Value = 290.8
$ws.Range("L2").Value = 938
$ws.Range("M2").Value = -177.8
$ws.Range("N2").Value = -1164
$ws.Range("H42").Value = 66998.5
$ws.Range("J42").Value = 66998.5
$ws.Range("L42").Value = 66998.5
$ws.Range("N42").Value = -67968.5
$ws.Range("H70").Value = 9611.875
$ws.Range("I70").Value = 9482.666999999999
$ws.Range("J70").Value = 9999.5
$ws.Range("K70").Value = 9482.666999999999
$ws.Range("L70").Value = 9999.5
$ws.Range("M70").Value = -9212.666999999999
$ws.Range("N70").Value = -10539.5
$ws.Range("H73").Value = 9611.875
$ws.Range("I73").Value = 9482.666999999999
$ws.Range("J73").Value = 9999.5
$ws.Range("K73").Value = 9482.666999999999
$ws.Range("L73").Value = 9999.5
$ws.Range("M73").Value = -8546.666999999999
$ws.Range("N73").Value = -11871.5
$ws.Range("H113").Value = 3466.3333
$ws.Range("I113").Value = 1199.5
$ws.Range("K113").Value = 1199.5
$ws.Range("M113").Value = 970.5
$ws.Range("H115").Value = 66998.5
$ws.Range("J115").Value = 66998.5
$ws.Range("L115").Value = 66998.5
$ws.Range("N115").Value = -69348.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3738.3076
$ws.Range("I16").Value = 4687.3335
$ws.Range("J16").Value = 2924.8572
$ws.Range("K16").Value = 4687.3335
$ws.Range("L16").Value = 2924.8572
$ws.Range("M16").Value = -4517.3335
$ws.Range("N16").Value = -3264.8572
$ws.Range("H82").Value = 7134.6924
$ws.Range("I82").Value = 7424.294
$ws.Range("J82").Value = 6587.6665
$ws.Range("K82").Value = 7424.294
$ws.Range("L82").Value = 6587.6665
$ws.Range("M82").Value = -7063.294
$ws.Range("N82").Value = -7309.6665
$ws.Range("H85").Value = 7134.6924
$ws.Range("I85").Value = 7424.294
$ws.Range("J85").Value = 6587.6665
$ws.Range("K85").Value = 7424.294
$ws.Range("L85").Value = 6587.6665
$ws.Range("M85").Value = -6176.294
$ws.Range("N85").Value = -9083.666499999999
$ws.Range("H87").Value = 56789
$ws.Range("J87").Value = 56789
$ws.Range("L87").Value = 56789
$ws.Range("N87").Value = -59035
$ws.Range("H88").Value = 39975
$ws.Range("J88").Value = 39975
$ws.Range("L88").Value = 39975
$ws.Range("N88").Value = -40831
$ws.Range("H90").Value = 56789
$ws.Range("J90").Value = 56789
$ws.Range("L90").Value = 170367
$ws.Range("N90").Value = -181599
$ws.Range("H91").Value = 39975
$ws.Range("J91").Value = 39975
$ws.Range("L91").Value = 39975
$ws.Range("N91").Value = -42939
$ws.Range("H136").Value = 6966.246
$ws.Range("I136").Value = 6527.2075
$ws.Range("K136").Value = 19581.6225
$ws.Range("M136").Value = -17031.6225

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14703.408
$ws.Range("I136").Value = 15532.021
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 46596.063
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -44046.063
$ws.Range("N136").Value = -11094
